# Apply crypto price/volume updates to match upstream scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.417.46'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '3.353.21'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '400.78'
$ws.Range('E5').Value = '  -3.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '125.75'
$ws.Range('E6').Value = '  +7.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.588'
$ws.Range('E7').Value = '  +2.14%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.657'
$ws.Range('E9').Value = '  +4.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.118'
$ws.Range('E10').Value = '  +1.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.90'
$ws.Range('E11').Value = '  +1.78%  '
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').Value = '3.891.79'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.25'
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.30'
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('D16').Value = '3.354.16'
$ws.Range('E16').Value = '  -1.32%  '
$ws.Range('D17').Value = '61.404.92'
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.16'
$ws.Range('E18').Value = '  +2.41%  '
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('E20').Value = '  +7.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '80.08'
$ws.Range('E22').Value = '  +7.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.61'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '298.74'
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.10'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('E26').Value = '  +11.14%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.93'
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.15'
$ws.Range('E28').Value = '  +5.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.50'
$ws.Range('E29').Value = '  -6.27%  '
$ws.Range('E30').Value = '  -2.13%  '
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.30'
$ws.Range('E33').Value = '  -1.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.49'
$ws.Range('E34').Value = '  -1.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '40.98'
$ws.Range('E35').Value = '  -5.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0477'
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.85'
$ws.Range('E37').Value = '  -0.96%  '
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('E39').Value = '  -1.29%  '
$ws.Range('E40').Value = '  -6.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '137.08'
$ws.Range('E41').Value = '  +2.74%  '
$ws.Range('E42').Value = '  +2.70%  '
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('E44').Value = '  -2.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.87'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.53'
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.05'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('D49').Value = '3.689.11'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').Value = '2.094.36'
$ws.Range('E50').Value = '  -3.60%  '
$ws.Range('E51').Value = '  -4.77%  '
